# AF-611 tests for named ranges and area intersections are moved to
# temporary_excel_files.
#
# The old D7/D8 helper cells (ISNA probes) are removed and D9 is turned
# into a plain SUM formula; the active selection follows the new last
# populated cell, D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two ISNA probe rows (D7: =ISNA(C5), D8: =ISNA(D6)).
$ws.Range("D7:D8").ClearContents()

# D9 becomes a simple SUM instead of =ISNA(D8).
$ws.Range("D9").Formula = "=SUM(A1,B1)"

# Match the saved selection/active cell in the workbook.
$ws.Range("D9").Select()
